$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "rose"
$ws.Range("B3").Value = 24.0
$ws.Range("C3").Value = "beijing"
$ws.Range("D3").Value = 12000.0
$ws.Range("E3").Value = $true

$ws.Range("A4").Value = "jack"
$ws.Range("B4").Value = 29.0
$ws.Range("C4").Value = "shanghai"
$ws.Range("D4").Value = 11111.1
$ws.Range("E4").Value = $true
